# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

# Rename header labels on the existing sheets
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Copy the header style (bold, bordered, centered) from the Weekly Quantity sheet
$wsWeekly.Range("A1:B1").Copy($ws3.Range("A1:B1"))
$wsWeekly.Range("A1:B1").Copy($ws3.Range("C1:D1"))

# Header row text (Copy above also copied values, so set the real text afterwards)
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Copy the date-format style from the Weekly Quantity sheet's "A" column onto the new "ds" column
$wsWeekly.Range("A2").Copy($ws3.Range("A2:A20"))

# Data rows
$ws3.Range("A2").Value = 45375.99999999999
$ws3.Range("B2").Value = 169
$ws3.Range("C2").Value = 73.14819228937866
$ws3.Range("D2").Value = 253.3032525501859
$ws3.Range("A3").Value = 45382.99999999999
$ws3.Range("B3").Value = 168
$ws3.Range("C3").Value = 77.83930190506764
$ws3.Range("D3").Value = 258.132958557022
$ws3.Range("A4").Value = 45396.99999999999
$ws3.Range("B4").Value = 164
$ws3.Range("C4").Value = 72.52580103907562
$ws3.Range("D4").Value = 253.3600910489346
$ws3.Range("A5").Value = 45403.99999999999
$ws3.Range("B5").Value = 163
$ws3.Range("C5").Value = 73.71888142736088
$ws3.Range("D5").Value = 244.3178572437851
$ws3.Range("A6").Value = 45417.99999999999
$ws3.Range("B6").Value = 160
$ws3.Range("C6").Value = 66.01941052962619
$ws3.Range("D6").Value = 237.4919050466811
$ws3.Range("A7").Value = 45424.99999999999
$ws3.Range("B7").Value = 158
$ws3.Range("C7").Value = 77.69488165789591
$ws3.Range("D7").Value = 245.3840534884889
$ws3.Range("A8").Value = 45431.99999999999
$ws3.Range("B8").Value = 157
$ws3.Range("C8").Value = 66.9290239867622
$ws3.Range("D8").Value = 247.1106353828449
$ws3.Range("A9").Value = 45438.99999999999
$ws3.Range("B9").Value = 155
$ws3.Range("C9").Value = 67.2293669439898
$ws3.Range("D9").Value = 247.4372618057503
$ws3.Range("A10").Value = 45445.99999999999
$ws3.Range("B10").Value = 154
$ws3.Range("C10").Value = 68.14330339638896
$ws3.Range("D10").Value = 241.7861403236214
$ws3.Range("A11").Value = 45501.99999999999
$ws3.Range("B11").Value = 142
$ws3.Range("C11").Value = 56.20435953048455
$ws3.Range("D11").Value = 227.5499951054729
$ws3.Range("A12").Value = 45508.99999999999
$ws3.Range("B12").Value = 140
$ws3.Range("C12").Value = 49.67985907524267
$ws3.Range("D12").Value = 231.0285029496172
$ws3.Range("A13").Value = 45515.99999999999
$ws3.Range("B13").Value = 139
$ws3.Range("C13").Value = 39.60418213823358
$ws3.Range("D13").Value = 222.8422938058911
$ws3.Range("A14").Value = 45522.99999999999
$ws3.Range("B14").Value = 137
$ws3.Range("C14").Value = 49.39888559743567
$ws3.Range("D14").Value = 219.7174475771813
$ws3.Range("A15").Value = 45529.99999999999
$ws3.Range("B15").Value = 135
$ws3.Range("C15").Value = 47.79347583919025
$ws3.Range("D15").Value = 217.4372493485614
$ws3.Range("A16").Value = 45536.99999999999
$ws3.Range("B16").Value = 134
$ws3.Range("C16").Value = 50.74364744473569
$ws3.Range("D16").Value = 212.5876040691496
$ws3.Range("A17").Value = 45543.99999999999
$ws3.Range("B17").Value = 132
$ws3.Range("C17").Value = 41.97252716677551
$ws3.Range("D17").Value = 217.7974705527842
$ws3.Range("A18").Value = 45550.99999999999
$ws3.Range("B18").Value = 131
$ws3.Range("C18").Value = 44.2122900199737
$ws3.Range("D18").Value = 211.7924529327219
$ws3.Range("A19").Value = 45557.99999999999
$ws3.Range("B19").Value = 129
$ws3.Range("C19").Value = 41.87259895342054
$ws3.Range("D19").Value = 219.4585154323543
$ws3.Range("A20").Value = 45564.99999999999
$ws3.Range("B20").Value = 128
$ws3.Range("C20").Value = 41.01698465937594
$ws3.Range("D20").Value = 215.5366743332275

$wsWeekly.Select()
